$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 header text changes from "Extra" to "Status"
$ws.Range("D1").Value = "Status"

# Selection moves from C2 to D2
$ws.Range("D2").Select()
